# Trade #26 closed at 2026-02-17 20:54:12 - unknown UNKNOWN +0.000%
#
# Applies:
#  - Summary sheet roll-up numbers
#  - Strategy Status row for MarketMaking
#  - All Trades: closes trade #54 (row 55) + appends new open trade #87 (row 88)
#  - MarketMaking: closes trade #54 (row 22) + appends new open trade #87 (row 55)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.55   # Current Capital
$summary.Range("B4").Value = 0.35      # Total P&L $
$summary.Range("B5").Value = 0.13      # Total P&L %
$summary.Range("B6").Value = 54        # Total Trades
$summary.Range("B7").Value = 26        # Winning Trades
$summary.Range("B9").Value = 48.15     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.55
$status.Range("D5").Value = 21
$status.Range("E5").Value = 0.24
$status.Range("F5").Value = 0.55
$status.Range("G5").Value = 57.14

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out existing trade #54 (row 55)
$allTrades.Range("G55").Value = 0.89
$allTrades.Range("H55").Value = "CLOSED"
$allTrades.Range("I55").Value = 4.7059
$allTrades.Range("J55").Value = 0.04
$allTrades.Range("K55").Value = 100.55
$allTrades.Range("L55").Value = "early_exit"
$allTrades.Range("M55").Value = 0.15

# Append newly opened trade #87 (row 88)
$allTrades.Range("A88").Value = 87
$allTrades.Range("B88").Value = "'2026-02-17"
$allTrades.Range("C88").Value = "20:54:05"
$allTrades.Range("D88").Value = "MarketMaking"
$allTrades.Range("E88").Value = "DOWN"
$allTrades.Range("F88").Value = 0.85
$allTrades.Range("H88").Value = "OPEN"
$allTrades.Range("I88").Value = 0
$allTrades.Range("J88").Value = 0
$allTrades.Range("K88").Value = 100.5134535840667
$allTrades.Range("M88").Value = 0
$allTrades.Range("N88").Value = 0
$allTrades.Range("O88").Value = 0
$allTrades.Range("P88").Value = 0.6
$allTrades.Range("Q88").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out existing trade #54 (row 22)
$mm.Range("G22").Value = 0.89
$mm.Range("H22").Value = "CLOSED"
$mm.Range("I22").Value = 4.7059
$mm.Range("J22").Value = 0.04
$mm.Range("K22").Value = 100.55
$mm.Range("P22").Value = "early_exit"
$mm.Range("Q22").Value = 0.15

# Append newly opened trade #87 (row 55)
$mm.Range("A55").Value = 87
$mm.Range("B55").Value = "'2026-02-17"
$mm.Range("C55").Value = "20:54:05"
$mm.Range("D55").Value = "MarketMaking"
$mm.Range("E55").Value = "DOWN"
$mm.Range("F55").Value = 0.85
$mm.Range("H55").Value = "OPEN"
$mm.Range("I55").Value = 0
$mm.Range("J55").Value = 0
$mm.Range("K55").Value = 100.5134535840667
$mm.Range("L55").Value = 0
$mm.Range("M55").Value = 0
$mm.Range("N55").Value = 0.6
$mm.Range("O55").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q55").Value = 0
